# 2021/11/16 - everyday data updated
# For each sheet, insert a new row at row 2 (pushing the existing dated rows
# down by one, carrying their formatting along) and populate the new row
# with the latest day's figures.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 台指期換倉成本計算 (A:日期 B:月份 C:結算價 D:未沖銷契約量 E:金額 F:成本) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2").Value = "日期：2021/11/16"
# "202112" looks purely numeric - force text so it isn't coerced to a number.
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "202112"
$ws1.Range("B2").ClearFormats()
$ws1.Range("C2").Value = 17680
$ws1.Range("D2").Value = 55055
$ws1.Range("E2").Value = 352344720
$ws1.Range("F2").Value = 17510

# --- Sheet 2: 散戶多空力道 (A:日期 B:散戶多空力道) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = "日期：2021/11/16"
$ws2.Range("B2").Value = 0.05

# --- Sheet 3: 三大法人買賣金額 (A:日期 B:外資 C:內資) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
$ws3.Range("A2").Value = "110年11月16日"
$ws3.Range("B2").Value = 76.95
$ws3.Range("C2").Value = -18.29

# --- Sheet 4: 大盤多空點位 (A:日期 B:隔日多空點位) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = "110年11月16日"
$ws4.Range("B2").Value = 17672.38

# --- Sheet 5: 期貨大額交易人未沖銷部位 (A:日期 B..N: ten/外資 figures) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
# "2021/11/16" looks like a date - force text so it isn't coerced to a date serial.
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "2021/11/16"
$ws5.Range("A2").ClearFormats()
$ws5.Range("B2").Value = 49536
$ws5.Range("C2").Value = 48555
$ws5.Range("D2").Value = 296
$ws5.Range("E2").Value = 81
$ws5.Range("F2").Value = 29676
$ws5.Range("G2").Value = 44739
$ws5.Range("H2").Value = 222
$ws5.Range("I2").Value = -294
$ws5.Range("J2").Value = -15063
$ws5.Range("K2").Value = 516
$ws5.Range("L2").Value = 74
$ws5.Range("M2").Value = 375
$ws5.Range("N2").Value = -301
